# Update automàtic: dades i banners [2026-02-10 19:50]
# Refresh of the meteo.cat daily summary extraction — updates the
# DATA_EXTRACCIO timestamps and the re-polled station readings.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-10 19:48:39"
$ws.Range("I2").Value = "43.0 mm"
$ws.Range("E3").Value = "2026-02-10 19:48:41"
$ws.Range("I3").Value = "24.8 mm"
$ws.Range("E4").Value = "2026-02-10 19:48:44"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "78%"
$ws.Range("J4").Value = "1003.9 hPa"
$ws.Range("O4").Value = "12.0 °C"
$ws.Range("E5").Value = "2026-02-10 19:48:46"
$ws.Range("I5").Value = "33.5 mm"
$ws.Range("E6").Value = "2026-02-10 19:48:48"
$ws.Range("J6").Value = "1004.2 hPa"
$ws.Range("O6").Value = "9.7 °C"
$ws.Range("E7").Value = "2026-02-10 19:48:51"
$ws.Range("E8").Value = "2026-02-10 19:48:53"
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = "84%"
$ws.Range("O8").Value = "12.0 °C"
$ws.Range("E9").Value = "2026-02-10 19:48:56"
$ws.Range("E10").Value = "2026-02-10 19:48:58"
$ws.Range("O10").Value = "10.2 °C"
$ws.Range("E11").Value = "2026-02-10 19:49:00"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = "89%"
$ws.Range("I11").Value = "2.3 mm"
$ws.Range("O11").Value = "7.7 °C"
$ws.Range("E12").Value = "2026-02-10 19:49:02"
$ws.Range("E13").Value = "2026-02-10 19:49:05"
$ws.Range("I13").Value = "7.0 mm"
$ws.Range("J13").Value = "1006.7 hPa"
$ws.Range("O13").Value = "5.2 °C"
$ws.Range("E14").Value = "2026-02-10 19:49:07"
$ws.Range("E15").Value = "2026-02-10 19:49:09"
$ws.Range("E16").Value = "2026-02-10 19:49:12"
$ws.Range("I16").Value = "25.1 mm"
$ws.Range("E17").Value = "2026-02-10 19:49:14"
$ws.Range("O17").Value = "4.7 °C"
$ws.Range("E18").Value = "2026-02-10 19:49:16"
$ws.Range("J18").Value = "1004.3 hPa"
$ws.Range("O18").Value = "10.3 °C"
$ws.Range("E19").Value = "2026-02-10 19:49:18"
$ws.Range("O19").Value = "6.6 °C"
$ws.Range("E20").Value = "2026-02-10 19:49:20"
$ws.Range("I20").Value = "8.1 mm"
$ws.Range("L20").Value = "76.7 km/h - 284º 19:09 TU"
$ws.Range("E21").Value = "2026-02-10 19:49:22"
$ws.Range("I21").Value = "8.7 mm"
$ws.Range("E22").Value = "2026-02-10 19:49:25"
$ws.Range("I22").Value = "9.2 mm"
$ws.Range("E23").Value = "2026-02-10 19:49:27"
$ws.Range("I23").Value = "25.6 mm"
$ws.Range("E24").Value = "2026-02-10 19:49:29"
$ws.Range("J24").Value = "1006.0 hPa"
$ws.Range("O24").Value = "11.3 °C"
$ws.Range("E25").Value = "2026-02-10 19:49:32"
$ws.Range("I25").Value = "18.5 mm"
$ws.Range("O25").Value = "1.6 °C"
$ws.Range("E26").Value = "2026-02-10 19:49:34"
$ws.Range("H26").NumberFormat = "@"
$ws.Range("H26").Value = "83%"
$ws.Range("I26").Value = "0.4 mm"
$ws.Range("E27").Value = "2026-02-10 19:49:37"
$ws.Range("I27").Value = "11.8 mm"
$ws.Range("O27").Value = "1.3 °C"
$ws.Range("E28").Value = "2026-02-10 19:49:39"
$ws.Range("H28").NumberFormat = "@"
$ws.Range("H28").Value = "82%"
$ws.Range("I28").Value = "0.3 mm"
$ws.Range("J28").Value = "1004.2 hPa"
$ws.Range("E29").Value = "2026-02-10 19:49:41"
$ws.Range("E30").Value = "2026-02-10 19:49:43"
$ws.Range("J30").Value = "1004.3 hPa"
$ws.Range("E31").Value = "2026-02-10 19:49:46"
$ws.Range("J31").Value = "1003.5 hPa"
$ws.Range("M31").Value = "14.5 °C 19:27 TU"
$ws.Range("O31").Value = "10.3 °C"
$ws.Range("E32").Value = "2026-02-10 19:49:48"
$ws.Range("H32").NumberFormat = "@"
$ws.Range("H32").Value = "90%"
$ws.Range("O32").Value = "10.5 °C"
$ws.Range("E33").Value = "2026-02-10 19:49:51"
$ws.Range("I33").Value = "10.7 mm"
$ws.Range("J33").Value = "1006.5 hPa"
$ws.Range("O33").Value = "4.2 °C"
$ws.Range("E34").Value = "2026-02-10 19:49:53"
$ws.Range("I34").Value = "13.3 mm"
$ws.Range("E35").Value = "2026-02-10 19:49:56"
$ws.Range("K35").Value = "10.4 MJ/m2"
$ws.Range("E36").Value = "2026-02-10 19:49:58"
$ws.Range("J36").Value = "1004.4 hPa"
$ws.Range("E37").Value = "2026-02-10 19:50:00"
$ws.Range("I37").Value = "0.9 mm"
$ws.Range("E38").Value = "2026-02-10 19:50:02"
$ws.Range("O38").Value = "10.9 °C"
$ws.Range("E39").Value = "2026-02-10 19:50:05"
$ws.Range("I39").Value = "10.1 mm"
$ws.Range("E40").Value = "2026-02-10 19:50:07"
$ws.Range("I40").Value = "13.2 mm"
$ws.Range("J40").Value = "1007.0 hPa"
$ws.Range("E41").Value = "2026-02-10 19:50:09"
$ws.Range("H41").NumberFormat = "@"
$ws.Range("H41").Value = "83%"
$ws.Range("L41").Value = "28.4 km/h - 228º 19:21 TU"
$ws.Range("E42").Value = "2026-02-10 19:50:12"
$ws.Range("E43").Value = "2026-02-10 19:50:14"
$ws.Range("O43").Value = "9.6 °C"
$ws.Range("E44").Value = "2026-02-10 19:50:16"
$ws.Range("I44").Value = "27.1 mm"
$ws.Range("E45").Value = "2026-02-10 19:50:18"
$ws.Range("I45").Value = "36.1 mm"
$ws.Range("E46").Value = "2026-02-10 19:50:21"
$ws.Range("H46").NumberFormat = "@"
$ws.Range("H46").Value = "82%"
$ws.Range("J46").Value = "1005.9 hPa"
$ws.Range("L46").Value = "49.0 km/h - 321º 19:20 TU"
$ws.Range("O46").Value = "14.3 °C"
